$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D9").Value = -8.417900000000005
$ws.Range("D13").Value = -7.351300000000002
$ws.Range("D16").Value = -8.477299999999996
$ws.Range("D18").Value = -8.388799999999993
$ws.Range("D20").Value = -8.0451
